# Update Lich Profits values per scheduled runner data refresh
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 3335.5881
$ws.Range("I92").Value = 3809.3333
$ws.Range("J92").Value = 2802.625
$ws.Range("K92").Value = 3809.3333
$ws.Range("L92").Value = 2802.625
$ws.Range("M92").Value = -2561.3333
$ws.Range("N92").Value = -5298.625
$ws.Range("H98").Value = 1431769.1
$ws.Range("J98").Value = 3339668.8
$ws.Range("L98").Value = 3339668.8
$ws.Range("N98").Value = -3342664.8
$ws.Range("H100").Value = 2616
$ws.Range("I100").Value = 2128.889
$ws.Range("J100").Value = 7000
$ws.Range("K100").Value = 2128.889
$ws.Range("L100").Value = 7000
$ws.Range("M100").Value = -1587.889
$ws.Range("N100").Value = -8082
$ws.Range("H106").Value = 11930
$ws.Range("I106").Value = 1884.6666
$ws.Range("K106").Value = 1884.6666
$ws.Range("M106").Value = -1253.6666
$ws.Range("H122").Value = 1431769.1
$ws.Range("J122").Value = 3339668.8
$ws.Range("L122").Value = 10019006.4
$ws.Range("N122").Value = -10023906.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 424.75
$ws.Range("I5").Value = 399.66666
$ws.Range("K5").Value = 399.66666
$ws.Range("M5").Value = -287.66666
$ws.Range("H32").Value = 10478.546
$ws.Range("I32").Value = 6171.9473
$ws.Range("K32").Value = 6171.9473
$ws.Range("M32").Value = -5884.9473
$ws.Range("H45").Value = 2852.6
$ws.Range("I45").Value = 2912.4443
$ws.Range("J45").Value = 2314
$ws.Range("K45").Value = 2912.4443
$ws.Range("L45").Value = 2314
$ws.Range("M45").Value = -2535.4443
$ws.Range("N45").Value = -3068
$ws.Range("H63").Value = 2089.4167
$ws.Range("J63").Value = 1899
$ws.Range("L63").Value = 1899
$ws.Range("N63").Value = -3271
$ws.Range("H66").Value = 2089.4167
$ws.Range("J66").Value = 1899
$ws.Range("L66").Value = 9495
$ws.Range("N66").Value = -16359
$ws.Range("H74").Value = 57805.5
$ws.Range("I74").Value = 66548.516
$ws.Range("K74").Value = 66548.516
$ws.Range("M74").Value = -65674.516
$ws.Range("H77").Value = 57805.5
$ws.Range("I77").Value = 66548.516
$ws.Range("K77").Value = 332742.58
$ws.Range("M77").Value = -328374.58
$ws.Range("H132").Value = 2127.9792
$ws.Range("I132").Value = 1697.25
$ws.Range("J132").Value = 2558.7083
$ws.Range("K132").Value = 5091.75
$ws.Range("L132").Value = 7676.124899999999
$ws.Range("M132").Value = -2561.75
$ws.Range("N132").Value = -12736.1249

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 424.75
$ws.Range("I4").Value = 399.66666
$ws.Range("K4").Value = 399.66666
$ws.Range("M4").Value = -284.66666
$ws.Range("H105").Value = 1758.7693
$ws.Range("I105").Value = 1545
$ws.Range("J105").Value = 2100.8
$ws.Range("K105").Value = 1545
$ws.Range("L105").Value = 2100.8
$ws.Range("M105").Value = 202
$ws.Range("N105").Value = -5594.8
$ws.Range("H134").Value = 4805.5864
$ws.Range("I134").Value = 4805.5864
$ws.Range("K134").Value = 14416.7592
$ws.Range("M134").Value = -11881.7592

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 627131.4399999999
$ws.Range("I31").Value = 771115.9399999999
$ws.Range("J31").Value = 3198.6667
$ws.Range("K31").Value = 771115.9399999999
$ws.Range("L31").Value = 3198.6667
$ws.Range("M31").Value = -770820.9399999999
$ws.Range("N31").Value = -3788.6667
$ws.Range("H34").Value = 627131.4399999999
$ws.Range("I34").Value = 771115.9399999999
$ws.Range("J34").Value = 3198.6667
$ws.Range("K34").Value = 771115.9399999999
$ws.Range("L34").Value = 3198.6667
$ws.Range("M34").Value = -770913.9399999999
$ws.Range("N34").Value = -3602.6667
$ws.Range("H58").Value = 2985.147
$ws.Range("I58").Value = 2821.6072
$ws.Range("J58").Value = 3748.3333
$ws.Range("K58").Value = 2821.6072
$ws.Range("L58").Value = 3748.3333
$ws.Range("M58").Value = -2618.6072
$ws.Range("N58").Value = -4154.3333
$ws.Range("H132").Value = 6269.8057
$ws.Range("I132").Value = 5495.857
$ws.Range("J132").Value = 8978.625
$ws.Range("K132").Value = 16487.571
$ws.Range("L132").Value = 26935.875
$ws.Range("M132").Value = -13957.571
$ws.Range("N132").Value = -31995.875
$ws.Range("H136").Value = 2985.147
$ws.Range("I136").Value = 2821.6072
$ws.Range("J136").Value = 3748.3333
$ws.Range("K136").Value = 8464.821599999999
$ws.Range("L136").Value = 11244.9999
$ws.Range("M136").Value = -5914.821599999999
$ws.Range("N136").Value = -16344.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 3647.8572
$ws.Range("J121").Value = 4048.625
$ws.Range("L121").Value = 12145.875
$ws.Range("N121").Value = -14765.875
$ws.Range("H131").Value = 1278.8055
$ws.Range("I131").Value = 623
$ws.Range("J131").Value = 1466.1786
$ws.Range("K131").Value = 1869
$ws.Range("L131").Value = 4398.5358
$ws.Range("M131").Value = 3171
$ws.Range("N131").Value = -14478.5358
$ws.Range("H140").Value = 10185.154
$ws.Range("I140").Value = 14635.934
$ws.Range("J140").Value = 4115.909
$ws.Range("K140").Value = 43907.802
$ws.Range("L140").Value = 12347.727
$ws.Range("M140").Value = -38727.802
$ws.Range("N140").Value = -22707.727

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 180.46666
$ws.Range("I2").Value = 186.21428
$ws.Range("K2").Value = 186.21428
$ws.Range("M2").Value = -73.21428
$ws.Range("H126").Value = 5388.636
$ws.Range("I126").Value = 6153.5557
$ws.Range("J126").Value = 1946.5
$ws.Range("K126").Value = 18460.6671
$ws.Range("L126").Value = 5839.5
$ws.Range("M126").Value = -15990.6671
$ws.Range("N126").Value = -10779.5
$ws.Range("H132").Value = 54949.8
$ws.Range("I132").Value = 60472.11
$ws.Range("K132").Value = 181416.33
$ws.Range("M132").Value = -178886.33

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4308.85
$ws.Range("I7").Value = 4068.6924
$ws.Range("J7").Value = 4754.857
$ws.Range("K7").Value = 4068.6924
$ws.Range("L7").Value = 4754.857
$ws.Range("M7").Value = -3956.6924
$ws.Range("N7").Value = -4978.857
$ws.Range("H22").Value = 5995.5
$ws.Range("J22").Value = 5995.5
$ws.Range("L22").Value = 5995.5
$ws.Range("N22").Value = -6585.5
$ws.Range("H27").Value = 5995.5
$ws.Range("J27").Value = 5995.5
$ws.Range("L27").Value = 5995.5
$ws.Range("N27").Value = -6209.5
$ws.Range("H46").Value = 2322.4285
$ws.Range("I46").Value = 1836.1305
$ws.Range("K46").Value = 1836.1305
$ws.Range("M46").Value = -1648.1305
$ws.Range("H82").Value = 3898.7144
$ws.Range("I82").Value = 3036
$ws.Range("J82").Value = 5049
$ws.Range("K82").Value = 3036
$ws.Range("L82").Value = 5049
$ws.Range("M82").Value = -2675
$ws.Range("N82").Value = -5771
$ws.Range("H85").Value = 3898.7144
$ws.Range("I85").Value = 3036
$ws.Range("J85").Value = 5049
$ws.Range("K85").Value = 3036
$ws.Range("L85").Value = 5049
$ws.Range("M85").Value = -1788
$ws.Range("N85").Value = -7545
$ws.Range("H126").Value = 4308.85
$ws.Range("I126").Value = 4068.6924
$ws.Range("J126").Value = 4754.857
$ws.Range("K126").Value = 12206.0772
$ws.Range("L126").Value = 14264.571
$ws.Range("M126").Value = -9736.0772
$ws.Range("N126").Value = -19204.571
$ws.Range("H132").Value = 1842.375
$ws.Range("I132").Value = 1834.1428
$ws.Range("K132").Value = 5502.428400000001
$ws.Range("M132").Value = -2972.428400000001
$ws.Range("H136").Value = 2528.3333
$ws.Range("I136").Value = 1793.3889
$ws.Range("K136").Value = 5380.1667
$ws.Range("M136").Value = -2830.1667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H103").Value = 15000
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()
$ws.Range("H132").Value = 7874.5
$ws.Range("I132").Value = 10197.6
$ws.Range("K132").Value = 30592.8
$ws.Range("M132").Value = -28062.8
$ws.Range("H136").Value = 386535.56
$ws.Range("I136").Value = 456469.8
$ws.Range("K136").Value = 1369409.4
$ws.Range("M136").Value = -1366859.4
$ws.Range("H141").Value = 64471.4
$ws.Range("J141").Value = 64471.4
$ws.Range("L141").Value = 64471.4
$ws.Range("N141").Value = -74831.39999999999

